# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3349
$ws1.Range("F4").Value = 61
$ws1.Range("F5").Value = 1443
$ws1.Range("F6").Value = 30
$ws1.Range("F7").Value = 319

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3349
$ws4.Range("F4").Value = 61
$ws4.Range("F5").Value = 1443
$ws4.Range("F6").Value = 30
$ws4.Range("F8").Value = 319
